$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Groups") updates
$ws.Range("C2").Value = 0.03134476197670545
$ws.Range("D2").Value = 0.03134476197670545
$ws.Range("E2").Value = 1.383680538000159
$ws.Range("G2").Value = 0.246

# Row 3 ("Residuals") updates
$ws.Range("C3").Value = 3.30736403593072
$ws.Range("D3").Value = 0.0226531783282926
